$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.078.10'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.830.82'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.81'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6242'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07491'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2921'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.12'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07680'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.830.67'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.011'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6669'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.71'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009352'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -6.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.975'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.079.71'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.077.31'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.56'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '222.52'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.125'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.28%  '
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.014'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.002'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.93'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1391'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.473'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.87'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.495'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05809'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +10.58%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.148'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.102'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.204'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7391'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.62%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.829'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.136'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.669'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.767'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.222.35'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.40%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01776'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.480'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.74%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8901'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.80'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000127'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.978.70'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.66'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5090'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('B50').Value = 'XinFinNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07584'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +13.37%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4054'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.34%  '
